$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 94 data
$ws.Range("A94").Value = "2024-01-15 18:58:58"
$ws.Range("B94").Value = 13
$ws.Range("C94").Value = 13
$ws.Range("D94").Value = 2
$ws.Range("E94").Value = 3
$ws.Range("F94").Value = 6
$ws.Range("G94").Value = 2
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0.001
$ws.Range("J94").Value = 0.05
$ws.Range("K94").Value = 0.003
$ws.Range("L94").Value = 100
$ws.Range("M94").Value = 500
$ws.Range("N94").Value = 10
$ws.Range("O94").Value = 5
$ws.Range("P94").Value = ""
$ws.Range("Q94").Value = "Data/bombay2.xlsx"

# New row 95 data
$ws.Range("A95").Value = "2024-01-16 13:32:44"
$ws.Range("B95").Value = 69
$ws.Range("C95").Value = 52
$ws.Range("D95").Value = 8
$ws.Range("E95").Value = 20
$ws.Range("F95").Value = 15
$ws.Range("G95").Value = 9
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0.001
$ws.Range("J95").Value = 0.05
$ws.Range("K95").Value = 0.003
$ws.Range("L95").Value = 100
$ws.Range("M95").Value = 500
$ws.Range("N95").Value = 10
$ws.Range("O95").Value = 5
$ws.Range("P95").Value = ""
$ws.Range("Q95").Value = "Data/bombay1.xlsx"
